$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("04")

# Rows 2-6 all hold the same invoice record (as shared strings). The new
# rows 7 and 8 are exact duplicates of that record, so copy an existing
# row (row 2) down into the two new rows rather than re-typing the values
# — this preserves the original shared-string cell typing (instead of
# Excel's usual "looks like a number" auto-conversion for values like
# "04", "601", "845.00") and keeps the default (unstyled) cell format.
$ws.Range("A2:N2").Copy()
$ws.Range("A7:N7").PasteSpecial()

$ws.Range("A2:N2").Copy()
$ws.Range("A8:N8").PasteSpecial()

$excel.CutCopyMode = $false
